# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp in the title row (A1)
# - Update El Salvador's active/recovered counts (row 73)
# - Give Georgia fresh covid numbers; since the sheet is sorted descending by
#   "Casos totales" (col B), Georgia's row now sits right after Jordania,
#   pushing Republica de Chipre / Liberia / Gambia down one row each (their
#   totals did not change, only their position)
# - Swap Islas Malvinas / Montserrat figures (tied total of 13, order swaps)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title timestamp
$ws.Range("A1").Value = "Datos actualizados a 10 de Agosto de 2020 a las 08:43"

# Row 73: El Salvador
$ws.Range("D73").Value = 9652
$ws.Range("E73").Value = 10222

# Row 144: Georgia (new figures, moved up in the sort order)
$ws.Range("A144").Value = "Georgia"
$ws.Range("B144").Value = 1250
$ws.Range("C144").Value = 25
$ws.Range("D144").Value = 1010
$ws.Range("E144").Value = 223
$ws.Range("H144").Value = 17

# Row 145: Republica de Chipre (shifted down, keeps its own old totals)
$ws.Range("A145").Value = "Republica de Chipre"
$ws.Range("B145").Value = 1242
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 870
$ws.Range("E145").Value = 353
$ws.Range("H145").Value = 19

# Row 146: Liberia (shifted down, keeps its own old totals)
$ws.Range("A146").Value = "Liberia"
$ws.Range("B146").Value = 1237
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 723
$ws.Range("E146").Value = 435
$ws.Range("H146").Value = 79

# Row 147: Gambia (shifted down, keeps its own old totals)
$ws.Range("A147").Value = "Gambia"
$ws.Range("B147").Value = 1235
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 221
$ws.Range("E147").Value = 991
$ws.Range("H147").Value = 23

# Row 202/203: Santa Lucia and Timor Oriental swap places (tied total = 25,
# identical figures across the board, only the country order changes)
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# Row 213/214: Islas Malvinas and Montserrat swap places (tied total = 13)
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
